{"js": "// Replace the date line and every \"NNN\u00f7N=\" division problem in the\n// practice table with the updated values from the new day's worksheet.\nconst replacements = [\n  [\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"],\n  [\"517\u00f72=\", \"797\u00f74=\"],\n  [\"735\u00f77=\", \"894\u00f78=\"],\n  [\"901\u00f77=\", \"978\u00f79=\"],\n  [\"961\u00f74=\", \"459\u00f76=\"],\n  [\"562\u00f73=\", \"190\u00f78=\"],\n  [\"994\u00f79=\", \"396\u00f72=\"],\n  [\"978\u00f76=\", \"107\u00f72=\"],\n  [\"141\u00f79=\", \"320\u00f79=\"],\n  [\"826\u00f75=\", \"980\u00f74=\"],\n  [\"273\u00f74=\", \"526\u00f75=\"],\n  [\"870\u00f79=\", \"622\u00f73=\"],\n  [\"139\u00f79=\", \"985\u00f73=\"],\n  [\"127\u00f72=\", \"319\u00f74=\"],\n  [\"165\u00f75=\", \"800\u00f76=\"],\n  [\"582\u00f76=\", \"816\u00f77=\"],\n  [\"306\u00f78=\", \"377\u00f73=\"],\n  [\"183\u00f76=\", \"456\u00f78=\"],\n  [\"666\u00f75=\", \"410\u00f73=\"],\n  [\"919\u00f73=\", \"626\u00f74=\"],\n  [\"856\u00f75=\", \"185\u00f78=\"],\n  [\"645\u00f76=\", \"848\u00f77=\"],\n  [\"752\u00f79=\", \"179\u00f76=\"],\n  [\"542\u00f73=\", \"386\u00f79=\"],\n  [\"919\u00f78=\", \"732\u00f74=\"],\n  [\"892\u00f76=\", \"884\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply targeted text replacements across the document (date line + all\n# \"NNN\u00f7N=\" division problems in the practice table), matching the diff.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1              # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\nReplace-Text \"2024-06-17 Monday\" \"2024-06-18 Tuesday\"\nReplace-Text \"517\u00f72=\" \"797\u00f74=\"\nReplace-Text \"735\u00f77=\" \"894\u00f78=\"\nReplace-Text \"901\u00f77=\" \"978\u00f79=\"\nReplace-Text \"961\u00f74=\" \"459\u00f76=\"\nReplace-Text \"562\u00f73=\" \"190\u00f78=\"\nReplace-Text \"994\u00f79=\" \"396\u00f72=\"\nReplace-Text \"978\u00f76=\" \"107\u00f72=\"\nReplace-Text \"141\u00f79=\" \"320\u00f79=\"\nReplace-Text \"826\u00f75=\" \"980\u00f74=\"\nReplace-Text \"273\u00f74=\" \"526\u00f75=\"\nReplace-Text \"870\u00f79=\" \"622\u00f73=\"\nReplace-Text \"139\u00f79=\" \"985\u00f73=\"\nReplace-Text \"127\u00f72=\" \"319\u00f74=\"\nReplace-Text \"165\u00f75=\" \"800\u00f76=\"\nReplace-Text \"582\u00f76=\" \"816\u00f77=\"\nReplace-Text \"306\u00f78=\" \"377\u00f73=\"\nReplace-Text \"183\u00f76=\" \"456\u00f78=\"\nReplace-Text \"666\u00f75=\" \"410\u00f73=\"\nReplace-Text \"919\u00f73=\" \"626\u00f74=\"\nReplace-Text \"856\u00f75=\" \"185\u00f78=\"\nReplace-Text \"645\u00f76=\" \"848\u00f77=\"\nReplace-Text \"752\u00f79=\" \"179\u00f76=\"\nReplace-Text \"542\u00f73=\" \"386\u00f79=\"\nReplace-Text \"919\u00f78=\" \"732\u00f74=\"\nReplace-Text \"892\u00f76=\" \"884\u00f76=\"\n"}
